$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 143095.72
$ws.Range("I9").Value = 200090
$ws.Range("J9").Value = 610
$ws.Range("K9").Value = 200090
$ws.Range("L9").Value = 610
$ws.Range("M9").Value = -199921
$ws.Range("N9").Value = -948
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H98").Value = 3924.5715
$ws.Range("I98").Value = 1242
$ws.Range("K98").Value = 1242
$ws.Range("M98").Value = 256
$ws.Range("H112").Value = 2568.8572
$ws.Range("J112").Value = 2717.25
$ws.Range("L112").Value = 8151.75
$ws.Range("N112").Value = -10367.75
$ws.Range("H122").Value = 3924.5715
$ws.Range("I122").Value = 1242
$ws.Range("K122").Value = 3726
$ws.Range("M122").Value = -1276
$ws.Range("H138").Value = 4702.3
$ws.Range("J138").Value = 4918.76
$ws.Range("L138").Value = 14756.28
$ws.Range("N138").Value = -25036.28

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9470.286
$ws.Range("I74").Value = 9470.286
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 9470.286
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -8596.286
$ws.Range("H77").Value = 9470.286
$ws.Range("I77").Value = 9470.286
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 47351.43
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -42983.43

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 67199.8
$ws.Range("J9").Value = 67199.8
$ws.Range("L9").Value = 67199.8
$ws.Range("N9").Value = -67535.8
$ws.Range("H86").Value = 7104.4546
$ws.Range("I86").Value = 2608.1667
$ws.Range("J86").Value = 12500
$ws.Range("K86").Value = 2608.1667
$ws.Range("L86").Value = 12500
$ws.Range("M86").Value = -1485.1667
$ws.Range("N86").Value = -14746
$ws.Range("H89").Value = 7104.4546
$ws.Range("I89").Value = 2608.1667
$ws.Range("J89").Value = 12500
$ws.Range("K89").Value = 13040.8335
$ws.Range("L89").Value = 62500
$ws.Range("M89").Value = -7424.833500000001
$ws.Range("N89").Value = -73732
$ws.Range("H107").Value = 1998.875
$ws.Range("I107").Value = 1832
$ws.Range("J107").Value = 2499.5
$ws.Range("K107").Value = 1832
$ws.Range("L107").Value = 2499.5
$ws.Range("M107").Value = 88
$ws.Range("N107").Value = -6339.5
$ws.Range("H135").Value = 68210.71000000001
$ws.Range("J135").Value = 68210.71000000001
$ws.Range("L135").Value = 68210.71000000001
$ws.Range("N135").Value = -78350.71000000001

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 11332.833
$ws.Range("I62").Value = 14999
$ws.Range("J62").Value = 7666.6665
$ws.Range("K62").Value = 14999
$ws.Range("L62").Value = 7666.6665
$ws.Range("M62").Value = -14375
$ws.Range("N62").Value = -8914.666499999999
$ws.Range("H65").Value = 11332.833
$ws.Range("I65").Value = 14999
$ws.Range("J65").Value = 7666.6665
$ws.Range("K65").Value = 74995
$ws.Range("L65").Value = 38333.3325
$ws.Range("M65").Value = -71875
$ws.Range("N65").Value = -44573.3325

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 2899.5
$ws.Range("J58").Value = 2899.5
$ws.Range("L58").Value = 8698.5
$ws.Range("N58").Value = -8954.5
$ws.Range("H122").Value = 1226.2858
$ws.Range("I122").Value = 546.5
$ws.Range("K122").Value = 4918.5
$ws.Range("M122").Value = -2468.5
$ws.Range("H123").Value = 4857.095
$ws.Range("I123").Value = 1999
$ws.Range("K123").Value = 5997
$ws.Range("M123").Value = -3547
$ws.Range("H132").Value = 14207.143
$ws.Range("J132").Value = 31416.666
$ws.Range("L132").Value = 282749.994
$ws.Range("N132").Value = -287809.994

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8479.799999999999
$ws.Range("I70").Value = 8479.799999999999
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 8479.799999999999
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -8209.799999999999
$ws.Range("H73").Value = 8479.799999999999
$ws.Range("I73").Value = 8479.799999999999
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 8479.799999999999
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -7543.799999999999
$ws.Range("H80").Value = 1830
$ws.Range("I80").Value = 2245
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 2245
$ws.Range("L80").Value = 1000
$ws.Range("M80").Value = -1247
$ws.Range("N80").Value = -2996
$ws.Range("H83").Value = 1830
$ws.Range("I83").Value = 2245
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 11225
$ws.Range("L83").Value = 5000
$ws.Range("M83").Value = -6233
$ws.Range("N83").Value = -14984

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 852.36365
$ws.Range("I55").Value = 595.25
$ws.Range("J55").Value = 999.2857
$ws.Range("K55").Value = 595.25
$ws.Range("L55").Value = 999.2857
$ws.Range("M55").Value = -422.25
$ws.Range("N55").Value = -1345.2857
$ws.Range("H61").Value = 1122.4445
$ws.Range("I61").Value = 871.7143
$ws.Range("K61").Value = 871.7143
$ws.Range("M61").Value = -669.7143
$ws.Range("H93").Value = 2311.6428
$ws.Range("I93").Value = 2082.4285
$ws.Range("J93").Value = 2540.8572
$ws.Range("K93").Value = 2082.4285
$ws.Range("L93").Value = 2540.8572
$ws.Range("M93").Value = -834.4285
$ws.Range("N93").Value = -5036.8572
$ws.Range("H100").Value = 1199.75
$ws.Range("I100").Value = 1199.75
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1199.75
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -658.75
$ws.Range("H113").Value = 1122.4445
$ws.Range("I113").Value = 871.7143
$ws.Range("K113").Value = 871.7143
$ws.Range("M113").Value = 1298.2857
$ws.Range("H136").Value = 3650.1667
$ws.Range("I136").Value = 2665.1667
$ws.Range("J136").Value = 4635.1665
$ws.Range("K136").Value = 7995.500100000001
$ws.Range("L136").Value = 13905.4995
$ws.Range("M136").Value = -5445.500100000001
$ws.Range("N136").Value = -19005.4995

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 750.7778
$ws.Range("I113").Value = 849.8333
$ws.Range("J113").Value = 701.25
$ws.Range("K113").Value = 2549.4999
$ws.Range("L113").Value = 2103.75
$ws.Range("M113").Value = -379.4998999999998
$ws.Range("N113").Value = -6443.75
$ws.Range("H132").Value = 4664.2
$ws.Range("I132").Value = 4465.8335
$ws.Range("J132").Value = 4961.75
$ws.Range("K132").Value = 13397.5005
$ws.Range("L132").Value = 14885.25
$ws.Range("M132").Value = -10867.5005
$ws.Range("N132").Value = -19945.25

Write-Host "All updates applied."